$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.278.48"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "1.982.02"
$ws.Range("E3").Value = "  +6.16%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'0.7857"
$ws.Range("E5").Value = "  +67.07%  "
$ws.Range("D6").Value = "'252.88"
$ws.Range("E6").Value = "  +3.75%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.3381"
$ws.Range("E8").Value = "  +17.62%  "
$ws.Range("D9").Value = "'25.59"
$ws.Range("E9").Value = "  +16.19%  "
$ws.Range("D10").Value = "'0.06927"
$ws.Range("E10").Value = "  +7.36%  "
$ws.Range("D11").Value = "'0.8319"
$ws.Range("E11").Value = "  +15.47%  "
$ws.Range("D12").Value = "'0.08107"
$ws.Range("E12").Value = "  +4.36%  "
$ws.Range("D13").Value = "1.988.68"
$ws.Range("E13").Value = "  +6.52%  "
$ws.Range("D14").Value = "'100.45"
$ws.Range("E14").Value = "  +4.81%  "
$ws.Range("D15").Value = "'5.437"
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("D16").Value = "'272.09"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "31.284.25"
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("D18").Value = "'13.85"
$ws.Range("E18").Value = "  +6.86%  "
$ws.Range("D19").Value = "'0.000007918"
$ws.Range("E19").Value = "  +5.21%  "
$ws.Range("D20").Value = "2.251.60"
$ws.Range("E20").Value = "  +6.45%  "
$ws.Range("D21").Value = "'5.707"
$ws.Range("E21").Value = "  +9.45%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'6.924"
$ws.Range("E24").Value = "  +11.40%  "
$ws.Range("D25").Value = "'9.614"
$ws.Range("E25").Value = "  +6.45%  "
$ws.Range("D26").Value = "'164.76"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "'0.1468"
$ws.Range("E27").Value = "  +53.10%  "
$ws.Range("D28").Value = "'19.72"
$ws.Range("E28").Value = "  +5.83%  "
$ws.Range("D29").Value = "'2.174"
$ws.Range("E29").Value = "  +16.11%  "
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("D31").Value = "'1.564"
$ws.Range("E31").Value = "  +6.72%  "
$ws.Range("D32").Value = "'4.542"
$ws.Range("E32").Value = "  +8.14%  "
$ws.Range("D33").Value = "'4.315"
$ws.Range("E33").Value = "  +5.38%  "
$ws.Range("D34").Value = "'0.05162"
$ws.Range("E34").Value = "  +7.44%  "
$ws.Range("D35").Value = "'1.208"
$ws.Range("E35").Value = "  +8.19%  "
$ws.Range("D36").Value = "'0.7513"
$ws.Range("E36").Value = "  +9.06%  "
$ws.Range("D37").Value = "'2.796"
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("D38").Value = "'1.002"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").Value = "'0.01998"
$ws.Range("E39").Value = "  +6.64%  "
$ws.Range("D40").Value = "'2.916"
$ws.Range("E40").Value = "  +3.84%  "
$ws.Range("D41").Value = "'6.615"
$ws.Range("E41").Value = "  +6.55%  "
$ws.Range("D42").Value = "'78.01"
$ws.Range("E42").Value = "  +5.29%  "
$ws.Range("D43").Value = "'0.4628"
$ws.Range("E43").Value = "  +9.91%  "
$ws.Range("D44").Value = "'2.055"
$ws.Range("E44").Value = "  +6.17%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.8522"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'105.29"
$ws.Range("E46").Value = "  +4.51%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "'9.954"
$ws.Range("E48").Value = "  +3.88%  "
$ws.Range("D49").Value = "'7.472"
$ws.Range("E49").Value = "  +7.67%  "
$ws.Range("D50").Value = "'36.40"
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("D51").Value = "'0.4264"
$ws.Range("E51").Value = "  +9.00%  "
